$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 38.25
$ws.Range("S2").Value = 59
$ws.Range("T2").Value = 65.5
$ws.Range("U2").Value = 127
$ws.Range("V2").Value = 55.82608695652174
$ws.Range("W2").Value = 25.09475761794446
$ws.Range("X2").Value = 52
$ws.Range("Y2").Value = 95.75
$ws.Range("Z2").Value = 124.5
$ws.Range("AA2").Value = 154.75
$ws.Range("AB2").Value = 231
$ws.Range("AC2").Value = 126.7173913043478
$ws.Range("AD2").Value = 43.6754001665714
$ws.Range("AE2").Value = 3
$ws.Range("AF2").Value = 6.25
$ws.Range("AG2").Value = 9
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 15
$ws.Range("AJ2").Value = 8.695652173913043
$ws.Range("AK2").Value = 3.133117476376037

# Row 5
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 26.25
$ws.Range("S5").Value = 40
$ws.Range("T5").Value = 54
$ws.Range("U5").Value = 90
$ws.Range("V5").Value = 41.91304347826087
$ws.Range("W5").Value = 19.66025442014254
$ws.Range("X5").Value = 36
$ws.Range("Y5").Value = 77.25
$ws.Range("Z5").Value = 105
$ws.Range("AA5").Value = 135
$ws.Range("AB5").Value = 174
$ws.Range("AC5").Value = 105.6304347826087
$ws.Range("AD5").Value = 38.93262481983868
$ws.Range("AE5").Value = 2
$ws.Range("AF5").Value = 6
$ws.Range("AG5").Value = 8
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 8.108695652173912
$ws.Range("AK5").Value = 3.07121590738236

# Row 6
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 19
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = 31
$ws.Range("U6").Value = 65
$ws.Range("V6").Value = 26.10869565217391
$ws.Range("W6").Value = 11.41972613189722
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 8
$ws.Range("Z6").Value = 9
$ws.Range("AA6").Value = 18.75
$ws.Range("AB6").Value = 43
$ws.Range("AC6").Value = 13.8695652173913
$ws.Range("AD6").Value = 11.34628416051739
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 1
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = 2
$ws.Range("AI6").Value = 5
$ws.Range("AJ6").Value = 1.521739130434783
$ws.Range("AK6").Value = 1.206263845005775

# Row 7
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = 21.25
$ws.Range("S7").Value = 30.5
$ws.Range("T7").Value = 45.75
$ws.Range("U7").Value = 108
$ws.Range("V7").Value = 34.39130434782609
$ws.Range("W7").Value = 19.09680399190697
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 27
$ws.Range("Z7").Value = 51.5
$ws.Range("AA7").Value = 73.25
$ws.Range("AB7").Value = 174
$ws.Range("AC7").Value = 53.3695652173913
$ws.Range("AD7").Value = 37.40847123054834
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 3
$ws.Range("AG7").Value = 4
$ws.Range("AH7").Value = 6.75
$ws.Range("AI7").Value = 11
$ws.Range("AJ7").Value = 4.434782608695652
$ws.Range("AK7").Value = 2.535544420286574

# Row 10
$ws.Range("Q10").Value = 6
$ws.Range("R10").Value = 16.25
$ws.Range("S10").Value = 23
$ws.Range("T10").Value = 32.75
$ws.Range("U10").Value = 49
$ws.Range("V10").Value = 25.26086956521739
$ws.Range("W10").Value = 10.84524428823517
$ws.Range("X10").Value = 0
$ws.Range("Y10").Value = 20
$ws.Range("Z10").Value = 37
$ws.Range("AA10").Value = 65.75
$ws.Range("AB10").Value = 141
$ws.Range("AC10").Value = 44.93478260869565
$ws.Range("AD10").Value = 33.15244396817227
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = 2
$ws.Range("AG10").Value = 3.5
$ws.Range("AH10").Value = 6
$ws.Range("AI10").Value = 10
$ws.Range("AJ10").Value = 4
$ws.Range("AK10").Value = 2.538591035287969

# Row 11
$ws.Range("Q11").Value = 11
$ws.Range("R11").Value = 22.25
$ws.Range("S11").Value = 29
$ws.Range("T11").Value = 38
$ws.Range("U11").Value = 112
$ws.Range("V11").Value = 33.91304347826087
$ws.Range("W11").Value = 20.72339106415906
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 8
$ws.Range("Z11").Value = 10
$ws.Range("AA11").Value = 19
$ws.Range("AB11").Value = 35
$ws.Range("AC11").Value = 12.52173913043478
$ws.Range("AD11").Value = 9.910127548085528
$ws.Range("AE11").Value = 0
$ws.Range("AF11").Value = 1
$ws.Range("AG11").Value = 1
$ws.Range("AH11").Value = 2
$ws.Range("AI11").Value = 4
$ws.Range("AJ11").Value = 1.347826086956522
$ws.Range("AK11").Value = 1.058665644507228

# Row 12
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 13
$ws.Range("S12").Value = 20
$ws.Range("T12").Value = 27.75
$ws.Range("U12").Value = 74
$ws.Range("V12").Value = 20.84074074074074
$ws.Range("W12").Value = 10.98565334310682
$ws.Range("X12").Value = 0
$ws.Range("Y12").Value = 22
$ws.Range("Z12").Value = 34
$ws.Range("AA12").Value = 44.75
$ws.Range("AB12").Value = 82
$ws.Range("AC12").Value = 34.0925925925926
$ws.Range("AD12").Value = 16.07587552584996
$ws.Range("AE12").Value = 0
$ws.Range("AF12").Value = 1
$ws.Range("AG12").Value = 2
$ws.Range("AH12").Value = 3
$ws.Range("AI12").Value = 5
$ws.Range("AJ12").Value = 2.192592592592593
$ws.Range("AK12").Value = 1.090658518696798

# Row 15
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 10
$ws.Range("S15").Value = 16
$ws.Range("T15").Value = 25
$ws.Range("U15").Value = 74
$ws.Range("V15").Value = 18.54074074074074
$ws.Range("W15").Value = 11.2869059003087
$ws.Range("X15").Value = 0
$ws.Range("Y15").Value = 25
$ws.Range("Z15").Value = 34
$ws.Range("AA15").Value = 47
$ws.Range("AB15").Value = 78
$ws.Range("AC15").Value = 35.62592592592593
$ws.Range("AD15").Value = 16.17023960143222
$ws.Range("AE15").Value = 0
$ws.Range("AF15").Value = 1
$ws.Range("AG15").Value = 2
$ws.Range("AH15").Value = 3
$ws.Range("AI15").Value = 5
$ws.Range("AJ15").Value = 2.240740740740741
$ws.Range("AK15").Value = 1.054908604106553

# Row 16
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = 9
$ws.Range("U16").Value = 42
$ws.Range("V16").Value = 7.018518518518518
$ws.Range("W16").Value = 6.563078944519698
$ws.Range("X16").Value = 0
$ws.Range("Y16").Value = 0
$ws.Range("Z16").Value = 8
$ws.Range("AA16").Value = 10
$ws.Range("AB16").Value = 34
$ws.Range("AC16").Value = 7.022222222222222
$ws.Range("AD16").Value = 8.471659876651998
$ws.Range("AE16").Value = 0
$ws.Range("AF16").Value = 0
$ws.Range("AG16").Value = 1
$ws.Range("AH16").Value = 1
$ws.Range("AI16").Value = 3
$ws.Range("AJ16").Value = 0.6851851851851852
$ws.Range("AK16").Value = 0.7813836282524065

# Row 17
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 14
$ws.Range("S17").Value = 24
$ws.Range("T17").Value = 32
$ws.Range("U17").Value = 57
$ws.Range("V17").Value = 23.61538461538462
$ws.Range("W17").Value = 11.60833193619993
$ws.Range("X17").Value = 0
$ws.Range("Y17").Value = 26
$ws.Range("Z17").Value = 35
$ws.Range("AA17").Value = 46.5
$ws.Range("AB17").Value = 78
$ws.Range("AC17").Value = 35.92307692307692
$ws.Range("AD17").Value = 14.22750460871335
$ws.Range("AE17").Value = 0
$ws.Range("AF17").Value = 2
$ws.Range("AG17").Value = 2
$ws.Range("AH17").Value = 3
$ws.Range("AI17").Value = 5
$ws.Range("AJ17").Value = 2.425641025641025
$ws.Range("AK17").Value = 0.9677047849863823

# Row 20
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 14.5
$ws.Range("S20").Value = 22
$ws.Range("T20").Value = 31
$ws.Range("U20").Value = 57
$ws.Range("V20").Value = 22.82051282051282
$ws.Range("W20").Value = 10.94708875971046
$ws.Range("X20").Value = 0
$ws.Range("Y20").Value = 24
$ws.Range("Z20").Value = 35
$ws.Range("AA20").Value = 45
$ws.Range("AB20").Value = 78
$ws.Range("AC20").Value = 34.97948717948718
$ws.Range("AD20").Value = 14.86986816071692
$ws.Range("AE20").Value = 0
$ws.Range("AF20").Value = 2
$ws.Range("AG20").Value = 2
$ws.Range("AH20").Value = 3
$ws.Range("AI20").Value = 6
$ws.Range("AJ20").Value = 2.394871794871795
$ws.Range("AK20").Value = 1.080747335442632

# Row 21
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 9
$ws.Range("T21").Value = 15
$ws.Range("U21").Value = 57
$ws.Range("V21").Value = 10.55384615384615
$ws.Range("W21").Value = 7.494391241408613
$ws.Range("X21").Value = 0
$ws.Range("Y21").Value = 0
$ws.Range("Z21").Value = 8
$ws.Range("AA21").Value = 16.5
$ws.Range("AB21").Value = 35
$ws.Range("AC21").Value = 9.353846153846154
$ws.Range("AD21").Value = 8.630616986278456
$ws.Range("AE21").Value = 0
$ws.Range("AF21").Value = 0
$ws.Range("AG21").Value = 1
$ws.Range("AH21").Value = 1.5
$ws.Range("AI21").Value = 3
$ws.Range("AJ21").Value = 0.9692307692307692
$ws.Range("AK21").Value = 0.8728888659183657
